$d = $word.ActiveDocument

$d.Content.Find.Execute("509÷8=63, 5", $true, $false, $false, $false, $false, $true, 1, $false, "628÷2=314, 0", 2) | Out-Null
$d.Content.Find.Execute("554÷8=69, 2", $true, $false, $false, $false, $false, $true, 1, $false, "302÷9=33, 5", 2) | Out-Null
$d.Content.Find.Execute("586÷3=195, 1", $true, $false, $false, $false, $false, $true, 1, $false, "780÷4=195, 0", 2) | Out-Null
$d.Content.Find.Execute("447÷6=74, 3", $true, $false, $false, $false, $false, $true, 1, $false, "727÷9=80, 7", 2) | Out-Null
$d.Content.Find.Execute("712÷4=178, 0", $true, $false, $false, $false, $false, $true, 1, $false, "197÷6=32, 5", 2) | Out-Null
$d.Content.Find.Execute("379÷8=47, 3", $true, $false, $false, $false, $false, $true, 1, $false, "394÷9=43, 7", 2) | Out-Null
$d.Content.Find.Execute("824÷4=206, 0", $true, $false, $false, $false, $false, $true, 1, $false, "177÷3=59, 0", 2) | Out-Null
$d.Content.Find.Execute("731÷5=146, 1", $true, $false, $false, $false, $false, $true, 1, $false, "931÷5=186, 1", 2) | Out-Null
$d.Content.Find.Execute("315÷5=63, 0", $true, $false, $false, $false, $false, $true, 1, $false, "509÷7=72, 5", 2) | Out-Null
$d.Content.Find.Execute("488÷2=244, 0", $true, $false, $false, $false, $false, $true, 1, $false, "398÷3=132, 2", 2) | Out-Null
$d.Content.Find.Execute("682÷6=113, 4", $true, $false, $false, $false, $false, $true, 1, $false, "566÷4=141, 2", 2) | Out-Null
$d.Content.Find.Execute("691÷8=86, 3", $true, $false, $false, $false, $false, $true, 1, $false, "114÷8=14, 2", 2) | Out-Null
$d.Content.Find.Execute("944÷5=188, 4", $true, $false, $false, $false, $false, $true, 1, $false, "990÷4=247, 2", 2) | Out-Null
$d.Content.Find.Execute("121÷9=13, 4", $true, $false, $false, $false, $false, $true, 1, $false, "757÷4=189, 1", 2) | Out-Null
$d.Content.Find.Execute("552÷5=110, 2", $true, $false, $false, $false, $false, $true, 1, $false, "407÷3=135, 2", 2) | Out-Null
$d.Content.Find.Execute("453÷3=151, 0", $true, $false, $false, $false, $false, $true, 1, $false, "521÷4=130, 1", 2) | Out-Null
$d.Content.Find.Execute("585÷3=195, 0", $true, $false, $false, $false, $false, $true, 1, $false, "477÷3=159, 0", 2) | Out-Null
$d.Content.Find.Execute("858÷7=122, 4", $true, $false, $false, $false, $false, $true, 1, $false, "514÷4=128, 2", 2) | Out-Null
$d.Content.Find.Execute("547÷3=182, 1", $true, $false, $false, $false, $false, $true, 1, $false, "510÷3=170, 0", 2) | Out-Null
$d.Content.Find.Execute("672÷4=168, 0", $true, $false, $false, $false, $false, $true, 1, $false, "173÷6=28, 5", 2) | Out-Null
$d.Content.Find.Execute("324÷6=54, 0", $true, $false, $false, $false, $false, $true, 1, $false, "755÷3=251, 2", 2) | Out-Null
$d.Content.Find.Execute("340÷5=68, 0", $true, $false, $false, $false, $false, $true, 1, $false, "230÷4=57, 2", 2) | Out-Null
$d.Content.Find.Execute("780÷2=390, 0", $true, $false, $false, $false, $false, $true, 1, $false, "648÷7=92, 4", 2) | Out-Null
$d.Content.Find.Execute("611÷6=101, 5", $true, $false, $false, $false, $false, $true, 1, $false, "137÷4=34, 1", 2) | Out-Null
$d.Content.Find.Execute("629÷2=314, 1", $true, $false, $false, $false, $false, $true, 1, $false, "270÷7=38, 4", 2) | Out-Null
